$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.57179999999999
$ws.Range("E4").Value = 16.3245
$ws.Range("C11").Value = -12.4878
$ws.Range("C12").Value = -12.04109999999999
$ws.Range("E14").Value = 16.73930000000001
$ws.Range("C15").Value = -13.5049
$ws.Range("E26").Value = 16.34199999999999
$ws.Range("C27").Value = -13.1464
$ws.Range("C28").Value = -13.60909999999999
$ws.Range("C31").Value = -13.4591
$ws.Range("E31").Value = 16.31080000000001
$ws.Range("C32").Value = -13.1849
$ws.Range("E35").Value = 16.4552
$ws.Range("C36").Value = -12.91230000000001
$ws.Range("E37").Value = 16.71780000000001
$ws.Range("C38").Value = -13.18990000000001
$ws.Range("E39").Value = 16.50199999999999
$ws.Range("E40").Value = 16.89820000000001
$ws.Range("E45").Value = 16.6469
$ws.Range("C46").Value = -14.65899999999999
$ws.Range("E52").Value = 17.1777
$ws.Range("C54").Value = -13.2802
$ws.Range("C55").Value = -13.76590000000001
$ws.Range("C56").Value = -12.61359999999999
$ws.Range("E57").Value = 16.75150000000001
$ws.Range("C67").Value = -10.61120000000001
$ws.Range("C69").Value = -11.97099999999999
$ws.Range("C72").Value = -11.6731
$ws.Range("C73").Value = -12.70290000000001
$ws.Range("E81").Value = 16.60779999999999
$ws.Range("C83").Value = -13.573
$ws.Range("E83").Value = 16.468
$ws.Range("C86").Value = -13.76569999999999
$ws.Range("C91").Value = -10.2605
$ws.Range("C93").Value = -10.99540000000001
$ws.Range("C99").Value = -13.8509
$ws.Range("E100").Value = 16.4248
$ws.Range("E102").Value = 16.8481
